$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price + volume) per diff.
# D-column (Price) values are forced to Text format before assignment
# so Excel does not auto-convert numeric-looking strings (e.g. "537.62")
# into actual numbers, matching the original inlineStr/text cell type.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.111.40"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  -0.46%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.521.63"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  +0.10%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  -0.16%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.62"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +0.34%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.51"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -1.65%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.567"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  +0.32%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.519.68"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -0.13%  "

# Row 10
$ws.Range("E10").Value = "  -0.10%  "

# Row 11
$ws.Range("E11").Value = "  -2.59%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.34"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -1.81%  "

# Row 13
$ws.Range("E13").Value = "  -2.99%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.967.01"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +0.04%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.09"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -2.15%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "59.020.35"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -0.47%  "

# Row 17
$ws.Range("E17").Value = "  -1.65%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.519.83"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -0.08%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.13"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -0.37%  "

# Row 20
$ws.Range("E20").Value = "  -0.76%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.34"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -0.34%  "

# Row 22
$ws.Range("E22").Value = "  -0.01%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.93"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +1.76%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.33"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +2.39%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.421"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -1.97%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.166"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -0.72%  "

# Row 27
$ws.Range("E27").Value = "  +0.14%  "

# Row 28
$ws.Range("E28").Value = "  -3.30%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.74"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -3.46%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0771"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -1.11%  "

# Row 31
$ws.Range("E31").Value = "  -0.69%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.78"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +0.20%  "

# Row 33
$ws.Range("E33").Value = "  +4.48%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.49"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +1.62%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -0.01%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.44"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -0.70%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.11"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -3.94%  "

# Row 38
$ws.Range("E38").Value = "  -3.04%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.64"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -0.73%  "

# Row 40
$ws.Range("E40").Value = "  -0.47%  "

# Row 41
$ws.Range("E41").Value = "  -2.05%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "286.06"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  +1.97%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.18"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -1.80%  "

# Row 44
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +0.16%  "

# Row 45
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "132.01"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +7.10%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.604"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +0.89%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.88"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +0.01%  "

# Row 48
$ws.Range("E48").Value = "  -0.33%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0508"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -1.18%  "

# Row 50
$ws.Range("E50").Value = "  -1.91%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.30"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -3.39%  "
